$p = $ppt.ActivePresentation
$s = $p.Slides.Item(27)
$s.HeadersFooters.Footer.Visible = $false
